$wb = $excel.ActiveWorkbook

# The "Repayment schedule" sheet (4th tab) gets a new blank column inserted
# before column N, shifting the existing "Late"/"Outstanding" columns (and
# the blank spacer column between them) one position to the right.
$ws4 = $wb.Worksheets.Item("Repayment schedule")
$ws4.Columns("N:N").Insert()

# The active/selected sheet moves from "Edit Repayment Schedule1" (tab 3)
# to "Repayment schedule" (tab 4), with a new selection on the latter.
$ws4.Activate()
$ws4.Range("S6").Select()
